$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$ws.Range("B2").Value = 0.1450617283950617
$ws.Range("C2").Value = 0.6666666666666666
$ws.Range("J2").Value = 0.0154320987654321
$ws.Range("P2").Value = 0.1141975308641975
$ws.Range("S2").Value = 0.05864197530864197
$ws.Range("B3").Value = 0.004484304932735426
$ws.Range("C3").Value = 0.03139013452914798
$ws.Range("J3").Value = 0.04484304932735426
$ws.Range("P3").Value = 0.7533632286995515
$ws.Range("S3").Value = 0.1659192825112108
$ws.Range("J4").Value = 0.03225806451612903
$ws.Range("P4").Value = 0.6774193548387096
$ws.Range("S4").Value = 0.2903225806451613
$ws.Range("B6").Value = 0.08962264150943396
$ws.Range("D6").Value = 0.01886792452830189
$ws.Range("F6").Value = 0.08490566037735849
$ws.Range("J6").Value = 0.2877358490566038
$ws.Range("O6").Value = 0.01415094339622642
$ws.Range("Q6").Value = 0.1132075471698113
$ws.Range("R6").Value = 0.07075471698113207
$ws.Range("S6").Value = 0.3207547169811321
$ws.Range("B7").Value = 0.1176470588235294
$ws.Range("D7").Value = 0.01176470588235294
$ws.Range("F7").Value = 0.02352941176470588
$ws.Range("J7").Value = 0.1235294117647059
$ws.Range("O7").Value = 0.04117647058823529
$ws.Range("Q7").Value = 0.1941176470588235
$ws.Range("R7").Value = 0.08823529411764706
$ws.Range("S7").Value = 0.4
$ws.Range("B8").Value = 0.1131221719457014
$ws.Range("D8").Value = 0.02714932126696833
$ws.Range("F8").Value = 0.04298642533936652
$ws.Range("J8").Value = 0.1221719457013575
$ws.Range("O8").Value = 0.01357466063348416
$ws.Range("Q8").Value = 0.2058823529411765
$ws.Range("R8").Value = 0.09954751131221719
$ws.Range("S8").Value = 0.3755656108597285
$ws.Range("B9").Value = 0.1354166666666667
$ws.Range("D9").Value = 0.03645833333333334
$ws.Range("F9").Value = 0.05208333333333334
$ws.Range("J9").Value = 0.06770833333333333
$ws.Range("O9").Value = 0.02083333333333333
$ws.Range("Q9").Value = 0.109375
$ws.Range("R9").Value = 0.1197916666666667
$ws.Range("S9").Value = 0.4583333333333333
$ws.Range("B10").Value = 0.1187077385424493
$ws.Range("D10").Value = 0.02779864763335838
$ws.Range("F10").Value = 0.06311044327573254
$ws.Range("J10").Value = 0.1359879789631856
$ws.Range("O10").Value = 0.01728024042073629
$ws.Range("Q10").Value = 0.2216378662659654
$ws.Range("R10").Value = 0.08489857250187828
$ws.Range("S10").Value = 0.3305785123966942
$ws.Range("G11").Value = 0.1550387596899225
$ws.Range("J11").Value = 0.06976744186046512
$ws.Range("K11").Value = 0.2054263565891473
$ws.Range("L11").Value = 0.5465116279069767
$ws.Range("S11").Value = 0.02325581395348837
$ws.Range("G12").Value = 0.7619047619047619
$ws.Range("J12").Value = 0.1836734693877551
$ws.Range("L12").Value = 0.0272108843537415
$ws.Range("S12").Value = 0.0272108843537415
$ws.Range("G13").Value = 0.6764705882352942
$ws.Range("J13").Value = 0.2352941176470588
$ws.Range("S13").Value = 0.08823529411764706
$ws.Range("F15").Value = 0.008064516129032258
$ws.Range("H15").Value = 0.1451612903225807
$ws.Range("I15").Value = 0.06854838709677419
$ws.Range("J15").Value = 0.375
$ws.Range("K15").Value = 0.06854838709677419
$ws.Range("M15").Value = 0.008064516129032258
$ws.Range("O15").Value = 0.09274193548387097
$ws.Range("S15").Value = 0.2338709677419355
$ws.Range("F16").Value = 0.02928870292887029
$ws.Range("H16").Value = 0.1841004184100418
$ws.Range("I16").Value = 0.06276150627615062
$ws.Range("J16").Value = 0.4058577405857741
$ws.Range("K16").Value = 0.100418410041841
$ws.Range("M16").Value = 0.03347280334728033
$ws.Range("O16").Value = 0.08368200836820083
$ws.Range("S16").Value = 0.100418410041841
$ws.Range("F17").Value = 0.02795698924731183
$ws.Range("H17").Value = 0.1698924731182796
$ws.Range("I17").Value = 0.09462365591397849
$ws.Range("J17").Value = 0.443010752688172
$ws.Range("K17").Value = 0.08172043010752689
$ws.Range("M17").Value = 0.01720430107526882
$ws.Range("O17").Value = 0.05806451612903226
$ws.Range("S17").Value = 0.1075268817204301
$ws.Range("F18").Value = 0.02380952380952381
$ws.Range("H18").Value = 0.2095238095238095
$ws.Range("I18").Value = 0.07142857142857142
$ws.Range("J18").Value = 0.4190476190476191
$ws.Range("K18").Value = 0.06666666666666667
$ws.Range("M18").Value = 0.01904761904761905
$ws.Range("O18").Value = 0.08571428571428572
$ws.Range("S18").Value = 0.1047619047619048
$ws.Range("F19").Value = 0.02161263507896924
$ws.Range("H19").Value = 0.200332502078138
$ws.Range("I19").Value = 0.08478802992518704
$ws.Range("J19").Value = 0.3840399002493766
$ws.Range("K19").Value = 0.08894430590191189
$ws.Range("M19").Value = 0.01080631753948462
$ws.Range("O19").Value = 0.07730673316708229
$ws.Range("S19").Value = 0.1321695760598504
